$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, which shifts the existing rows 308-341
# down to 309-342 (data for those rows is preserved unchanged).
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row 308 with the new weekly price record.
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = 44946
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = 100112032
$ws.Range("G308").Value = "Zapallo italiano"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 300
$ws.Range("K308").Value = 14000
$ws.Range("L308").Value = 15000
$ws.Range("M308").Value = 14500
$ws.Range("N308").Value = "$/caja 50 unidades"
$ws.Range("O308").Value = "Región de O'Higgins"
$ws.Range("P308").Value = 290
$ws.Range("Q308").Value = 50
$ws.Range("R308").Value = "Hortaliza"

# Match the date cell's number format (a date/time format) used by the
# other rows' "Fecha" column.
$ws.Range("D308").NumberFormat = $ws.Range("D309").NumberFormat
